$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Replace the word "nobildētos" with "noskenētos" in the task description in B5
$ws.Range("B5").Value = "Izveidot sistēmu, ar ko var vērtēt pārbaudes darbus, izmantojot to pareizo atbilžu izkārtojumus un noskenētos pārbaudes darbus"

# Update the active selection to match the saved workbook state
$ws.Range("L8").Select()
